$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Property" -> "DataNode"
# (commit: "unify the conception of DataNode, DataTable, Entity")
$ws.Name = "DataNode"

# Nudge the stored column widths a hair closer to the values the
# workbook carries after being re-saved (23.875 / 8.125 chars); the
# underlying character counts are effectively unchanged, we're just
# re-asserting them so they get persisted with this engine's own
# (Windows-style) width formula instead of the original Mac-Excel one.
$ws.Columns.Item(1).ColumnWidth = 23.142857142857142
$ws.Columns.Item(2).ColumnWidth = 7.428571428571429

# Leave the active selection where the author last left it.
$ws.Range("D39").Select() | Out-Null
